$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 590.65625
$ws.Range("J17").Value = 590.65625
$ws.Range("L17").Value = 1771.96875
$ws.Range("N17").Value = -2107.96875
$ws.Range("H28").Value = 2356.2144
$ws.Range("I28").Value = 1191.8889
$ws.Range("K28").Value = 1191.8889
$ws.Range("M28").Value = -706.8888999999999
$ws.Range("H62").Value = 4510.125
$ws.Range("I62").Value = 3394.5
$ws.Range("J62").Value = 4882
$ws.Range("K62").Value = 3394.5
$ws.Range("L62").Value = 4882
$ws.Range("M62").Value = -2770.5
$ws.Range("N62").Value = -6130
$ws.Range("H65").Value = 4510.125
$ws.Range("I65").Value = 3394.5
$ws.Range("J65").Value = 4882
$ws.Range("K65").Value = 16972.5
$ws.Range("L65").Value = 24410
$ws.Range("M65").Value = -13852.5
$ws.Range("N65").Value = -30650
$ws.Range("H98").Value = 3043.3333
$ws.Range("I98").Value = 3043.3333
$ws.Range("K98").Value = 3043.3333
$ws.Range("M98").Value = -1545.3333
$ws.Range("H112").Value = 4243.6484
$ws.Range("J112").Value = 4243.6484
$ws.Range("L112").Value = 12730.9452
$ws.Range("N112").Value = -14946.9452
$ws.Range("H122").Value = 3043.3333
$ws.Range("I122").Value = 3043.3333
$ws.Range("K122").Value = 9129.999899999999
$ws.Range("M122").Value = -6679.999899999999
$ws.Range("H133").Value = 79101.53999999999
$ws.Range("J133").Value = 79101.53999999999
$ws.Range("L133").Value = 79101.53999999999
$ws.Range("N133").Value = -89221.53999999999
$ws.Range("H135").Value = 2861.7097
$ws.Range("I135").Value = 1595.2354
$ws.Range("J135").Value = 4399.5713
$ws.Range("K135").Value = 14357.1186
$ws.Range("L135").Value = 39596.14169999999
$ws.Range("M135").Value = -11822.1186
$ws.Range("N135").Value = -44666.14169999999
$ws.Range("H138").Value = 8257.938
$ws.Range("J138").Value = 8991.721
$ws.Range("L138").Value = 26975.163
$ws.Range("N138").Value = -37255.163

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2351.3635
$ws.Range("I45").Value = 1239.5625
$ws.Range("K45").Value = 1239.5625
$ws.Range("M45").Value = -862.5625
$ws.Range("H61").Value = 14058.923
$ws.Range("I61").Value = 12995.667
$ws.Range("K61").Value = 12995.667
$ws.Range("M61").Value = -12783.667
$ws.Range("H80").Value = 89957.39999999999
$ws.Range("J80").Value = 89957.39999999999
$ws.Range("L80").Value = 89957.39999999999
$ws.Range("N80").Value = -91953.39999999999
$ws.Range("H83").Value = 89957.39999999999
$ws.Range("J83").Value = 89957.39999999999
$ws.Range("L83").Value = 269872.2
$ws.Range("N83").Value = -279856.2
$ws.Range("H102").Value = 807062.4399999999
$ws.Range("I102").Value = 1246592
$ws.Range("J102").Value = 1258.1666
$ws.Range("K102").Value = 1246592
$ws.Range("L102").Value = 1258.1666
$ws.Range("M102").Value = -1244970
$ws.Range("N102").Value = -4502.1666
$ws.Range("H136").Value = 14058.923
$ws.Range("I136").Value = 12995.667
$ws.Range("K136").Value = 38987.001
$ws.Range("M136").Value = -36437.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 35768
$ws.Range("I82").Value = 9244.5
$ws.Range("J82").Value = 56986.8
$ws.Range("K82").Value = 9244.5
$ws.Range("L82").Value = 56986.8
$ws.Range("M82").Value = -8861.5
$ws.Range("N82").Value = -57752.8
$ws.Range("H85").Value = 35768
$ws.Range("I85").Value = 9244.5
$ws.Range("J85").Value = 56986.8
$ws.Range("K85").Value = 9244.5
$ws.Range("L85").Value = 56986.8
$ws.Range("M85").Value = -7918.5
$ws.Range("N85").Value = -59638.8
$ws.Range("H135").Value = 85186.664
$ws.Range("J135").Value = 85186.664
$ws.Range("L135").Value = 85186.664
$ws.Range("N135").Value = -95326.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1431758.1
$ws.Range("I58").Value = 1669134.9
$ws.Range("K58").Value = 1669134.9
$ws.Range("M58").Value = -1668931.9
$ws.Range("H60").Value = 99499.75
$ws.Range("J60").Value = 99499.75
$ws.Range("L60").Value = 99499.75
$ws.Range("N60").Value = -100521.75
$ws.Range("H68").Value = 70018.78999999999
$ws.Range("J68").Value = 71558.69500000001
$ws.Range("L68").Value = 71558.69500000001
$ws.Range("N68").Value = -73056.69500000001
$ws.Range("H71").Value = 70018.78999999999
$ws.Range("J71").Value = 71558.69500000001
$ws.Range("L71").Value = 214676.085
$ws.Range("N71").Value = -222164.085
$ws.Range("H74").Value = 62481.637
$ws.Range("J74").Value = 62481.637
$ws.Range("L74").Value = 62481.637
$ws.Range("N74").Value = -64229.637
$ws.Range("H77").Value = 62481.637
$ws.Range("J77").Value = 62481.637
$ws.Range("L77").Value = 187444.911
$ws.Range("N77").Value = -196180.911
$ws.Range("H136").Value = 1431758.1
$ws.Range("I136").Value = 1669134.9
$ws.Range("K136").Value = 5007404.699999999
$ws.Range("M136").Value = -5004854.699999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1628803.5
$ws.Range("I131").Value = 1943.7858
$ws.Range("J131").Value = 2171090
$ws.Range("K131").Value = 5831.357400000001
$ws.Range("L131").Value = 6513270
$ws.Range("M131").Value = -791.3574000000008
$ws.Range("N131").Value = -6523350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 21742
$ws.Range("H57").Value = 11999
$ws.Range("I57").Value = 11999
$ws.Range("K57").Value = 11999
$ws.Range("M57").Value = -11179
$ws.Range("H122").Value = 530041.5600000001
$ws.Range("I122").Value = 692904.7
$ws.Range("K122").Value = 2078714.1
$ws.Range("M122").Value = -2076264.1
$ws.Range("H126").Value = 4641.0586
$ws.Range("I126").Value = 2212.5
$ws.Range("K126").Value = 6637.5
$ws.Range("M126").Value = -4167.5
$ws.Range("H132").Value = 3146.1052
$ws.Range("I132").Value = 3072.394
$ws.Range("K132").Value = 9217.181999999999
$ws.Range("M132").Value = -6687.181999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 12504729
$ws.Range("I40").Value = 25002962
$ws.Range("J40").Value = 6496.5
$ws.Range("K40").Value = 25002962
$ws.Range("L40").Value = 6496.5
$ws.Range("M40").Value = -25002826
$ws.Range("N40").Value = -6768.5
$ws.Range("H46").Value = 5515.5415
$ws.Range("J46").Value = 5705.857
$ws.Range("L46").Value = 5705.857
$ws.Range("N46").Value = -6081.857
$ws.Range("H61").Value = 3945.0667
$ws.Range("I61").Value = 3922.24
$ws.Range("K61").Value = 3922.24
$ws.Range("M61").Value = -3720.24
$ws.Range("H113").Value = 3945.0667
$ws.Range("I113").Value = 3922.24
$ws.Range("K113").Value = 3922.24
$ws.Range("M113").Value = -1752.24
$ws.Range("H122").Value = 8000
$ws.Range("J122").Value = 8000
$ws.Range("L122").Value = 24000
$ws.Range("N122").Value = -28900
$ws.Range("H132").Value = 2834.33
$ws.Range("I132").Value = 2839.7937
$ws.Range("K132").Value = 8519.381100000001
$ws.Range("M132").Value = -5989.381100000001
$ws.Range("H136").Value = 4264.3
$ws.Range("I136").Value = 4004.5774
$ws.Range("J136").Value = 4900.1724
$ws.Range("K136").Value = 12013.7322
$ws.Range("L136").Value = 14700.5172
$ws.Range("M136").Value = -9463.7322
$ws.Range("N136").Value = -19800.5172

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 203743.69
$ws.Range("I14").Value = 495316.66
$ws.Range("J14").Value = 28799.9
$ws.Range("K14").Value = 495316.66
$ws.Range("L14").Value = 28799.9
$ws.Range("M14").Value = -495148.66
$ws.Range("N14").Value = -29135.9
$ws.Range("H64").Value = 71328.25
$ws.Range("I64").Value = 44990
$ws.Range("K64").Value = 44990
$ws.Range("M64").Value = -44742
$ws.Range("H67").Value = 71328.25
$ws.Range("I67").Value = 44990
$ws.Range("K67").Value = 44990
$ws.Range("M67").Value = -44132
$ws.Range("H100").Value = 801258.9399999999
$ws.Range("I100").Value = 1053382.9
$ws.Range("K100").Value = 2106765.8
$ws.Range("M100").Value = -2106224.8
$ws.Range("H122").Value = 3330.913
$ws.Range("I122").Value = 2563.7368
$ws.Range("K122").Value = 7691.2104
$ws.Range("M122").Value = -5241.2104
